# This workbook is a daily "accident report" export. The edit replaces the
# previous day's rows of accident data (2018-11-12, row 2 through row 21)
# with a new day's rows (2018-11-14, row 2 through row 20) - i.e. one fewer
# data row overall, plus a newly introduced "No Injuries" category value
# alongside the pre-existing "Injuries" / "Unknown Injuries" values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-20 (columns A, C, D, E, F, G, H, K, P - all other
# columns in these rows stay blank, same as before the edit). K=$null means
# the City cell should end up empty for that row.
$rows = @(
    @{ Row=2;  A=3;   C="Unknown Injuries"; D=35.015818; E=-85.254192;  F="2018-11-14"; G="23:00:59"; H="18300 Interstate 24 Eb";          K="EAST RIDGE";      P="23" },
    @{ Row=3;  A=7;   C="Injuries";         D=35.036773; E=-85.264846;  F="2018-11-14"; G="22:53:16"; H="2324 E 4th St";                   K="CHATTANOOGA";     P="22" },
    @{ Row=4;  A=8;   C="Injuries";         D=35.036773; E=-85.264846;  F="2018-11-14"; G="22:53:03"; H="2324 E 4th St";                   K="CHATTANOOGA";     P="22" },
    @{ Row=5;  A=10;  C="Unknown Injuries"; D=35.002679; E=-85.21117;   F="2018-11-14"; G="22:35:46"; H="130 INTERSTATE 75 NB";             K="CHATTANOOGA";     P="22" },
    @{ Row=6;  A=14;  C="Injuries";         D=35.019215; E=-85.298212;  F="2018-11-14"; G="21:05:14"; H="2707 Cannon Ave";                  K="CHATTANOOGA";     P="21" },
    @{ Row=7;  A=16;  C="Injuries";         D=35.019215; E=-85.298212;  F="2018-11-14"; G="21:04:34"; H="2707 Cannon Ave";                  K="CHATTANOOGA";     P="21" },
    @{ Row=8;  A=23;  C="Unknown Injuries"; D=34.996081; E=-85.206737;  F="2018-11-14"; G="19:24:37"; H="100 Interstate 75 Nb";             K=$null;             P="19" },
    @{ Row=9;  A=39;  C="Injuries";         D=35.081593; E=-85.209734;  F="2018-11-14"; G="17:48:31"; H="Bonny Oaks Dr / Highway 58";       K="CHATTANOOGA";     P="17" },
    @{ Row=10; A=40;  C="Unknown Injuries"; D=35.081593; E=-85.209734;  F="2018-11-14"; G="17:47:42"; H="Highway 58 / Bonny Oaks Dr";       K="CHATTANOOGA";     P="17" },
    @{ Row=11; A=40;  C="Unknown Injuries"; D=35.081593; E=-85.209734;  F="2018-11-14"; G="17:47:42"; H="Highway 58 / Bonny Oaks Dr";       K="CHATTANOOGA";     P="17" },
    @{ Row=12; A=48;  C="Injuries";         D=35.127027; E=-85.029316;  F="2018-11-14"; G="17:12:05"; H="1540 Interstate 75 Sb";            K="HAMILTON COUNTY"; P="17" },
    @{ Row=13; A=62;  C="Injuries";         D=35.019438; E=-85.283738;  F="2018-11-14"; G="16:09:56"; H="S Orchard Knob Ave / E 23rd St";   K="CHATTANOOGA";     P="16" },
    @{ Row=14; A=76;  C="Injuries";         D=35.043072; E=-85.279331;  F="2018-11-14"; G="15:11:31"; H="1330 E 3rd St";                    K="CHATTANOOGA";     P="15" },
    @{ Row=15; A=84;  C="No Injuries";      D=35.01068;  E=-85.244059;  F="2018-11-14"; G="14:27:10"; H="BELVOIR AVE / S TERRACE";          K="CHATTANOOGA";     P="14" },
    @{ Row=16; A=85;  C="Unknown Injuries"; D=35.01068;  E=-85.244059;  F="2018-11-14"; G="14:27:00"; H="S TERRACE / BELVOIR AVE";          K="EAST RIDGE";      P="14" },
    @{ Row=17; A=92;  C="Injuries";         D=35.039678; E=-85.27432;   F="2018-11-14"; G="13:22:53"; H="300 N Hawthorne St";               K="CHATTANOOGA";     P="13" },
    @{ Row=18; A=93;  C="Injuries";         D=35.039678; E=-85.27432;   F="2018-11-14"; G="13:22:40"; H="300 N HAWTHORNE ST";               K="CHATTANOOGA";     P="13" },
    @{ Row=19; A=110; C="Unknown Injuries"; D=35.066856; E=-85.244523;  F="2018-11-14"; G="11:40:15"; H="2720 GLASS ST";                    K="CHATTANOOGA";     P="11" },
    @{ Row=20; A=128; C="Injuries";         D=35.026626; E=-85.310189;  F="2018-11-14"; G="09:07:02"; H="2500 Market St";                  K=$null;             P="9"  }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $r.A    # A - Accident
    $ws.Cells.Item($row, 3).Value  = $r.C    # C - Problem
    $ws.Cells.Item($row, 4).Value  = $r.D    # D - Latitude
    $ws.Cells.Item($row, 5).Value  = $r.E    # E - Longitude
    $ws.Cells.Item($row, 6).Value  = $r.F    # F - Date
    $ws.Cells.Item($row, 7).Value  = $r.G    # G - Time
    $ws.Cells.Item($row, 8).Value  = $r.H    # H - Address

    if ($null -eq $r.K) {
        $ws.Cells.Item($row, 11).ClearContents()   # K - City (now blank)
    } else {
        $ws.Cells.Item($row, 11).Value = $r.K       # K - City
    }

    $ws.Cells.Item($row, 16).Value = $r.P    # P - Hour
}

# Row 21 (accident 125, SODDY DAISY / HAMILTON COUNTY "Entrapment" record)
# no longer exists in the new data - remove it entirely, shifting the
# worksheet dimension from A1:BA21 down to A1:BA20.
$ws.Rows(21).Delete()
